$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "303.47"
Set-TextValue "E2" "4.58%"
Set-TextValue "D3" "32.13"
Set-TextValue "E3" "9.71%"
Set-TextValue "D4" "5.253"
Set-TextValue "E4" "-0.59%"
Set-TextValue "D5" "0.07523"
Set-TextValue "E5" "5.16%"
Set-TextValue "D6" "7.894"
Set-TextValue "E6" "5.69%"
Set-TextValue "D7" "3.814"
Set-TextValue "E7" "6.81%"
Set-TextValue "D8" "1.511"
Set-TextValue "E8" "7.56%"
Set-TextValue "D9" "0.9222"
Set-TextValue "D10" "0.1699"
Set-TextValue "E10" "4.90%"
Set-TextValue "D11" "0.08002"
Set-TextValue "E11" "4.58%"
Set-TextValue "D12" "0.08022"
Set-TextValue "E12" "3.35%"
Set-TextValue "D13" "0.03039"
Set-TextValue "E13" "4.02%"
Set-TextValue "D14" "0.09898"
Set-TextValue "E14" "9.71%"
Set-TextValue "D15" "0.001492"
Set-TextValue "E15" "-6.28%"
Set-TextValue "D16" "0.04592"
Set-TextValue "D17" "0.006582"
Set-TextValue "E17" "6.19%"
Set-TextValue "D18" "3.459"
Set-TextValue "E18" "-0.81%"
Set-TextValue "D19" "2.231"
Set-TextValue "E19" "-0.05%"
Set-TextValue "D20" "0.3303"
Set-TextValue "E20" "2.02%"
Set-TextValue "E21" "-0.45%"
Set-TextValue "D22" "4.497"
Set-TextValue "E22" "12.05%"
Set-TextValue "D23" "0.1620"
Set-TextValue "E23" "1.38%"
Set-TextValue "D24" "0.001216"
Set-TextValue "E24" "0.61%"
Set-TextValue "D25" "0.004455"
Set-TextValue "E25" "5.99%"
Set-TextValue "D26" "0.0001397"
Set-TextValue "E26" "19.58%"
Set-TextValue "D27" "0.0001778"
Set-TextValue "E27" "6.26%"
Set-TextValue "D39" "0.01711"
Set-TextValue "E39" "2,525.14%"
Set-TextValue "D40" "0.04484"
Set-TextValue "E40" "1.71%"
Set-TextValue "D41" "0.006967"
Set-TextValue "E41" "-0.76%"
Set-TextValue "D42" "0.1350"
Set-TextValue "E42" "6.91%"
Set-TextValue "D43" "0.002136"
Set-TextValue "E43" "3.30%"
Set-TextValue "D44" "0.01285"
Set-TextValue "E44" "9.08%"
Set-TextValue "D45" "0.00006154"
Set-TextValue "E45" "5.26%"
Set-TextValue "D46" "1.862"
Set-TextValue "E46" "-3.49%"
Set-TextValue "D47" "0.01497"
Set-TextValue "E47" "15.39%"
